$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 163.680235902
$ws.Range("C2").Value = 0.00000006913746076069999

$ws.Range("B3").Value = 163680.235902
$ws.Range("C3").Value = 0.0006913746076069999

$ws.Range("B4").Value = 327360.471804
$ws.Range("C4").Value = 0.002074123822820999

$ws.Range("B5").Value = 6547.209436079999
$ws.Range("C5").Value = 0.00004148247645641999
